$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("C13").Value = ""

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "NSE:5PAISA"
$ws.Range("C2").Value = "NSE:BIGBLOC"
$ws.Range("E2").Value = "NSE:GODREJCP"
$ws.Range("F2").Value = "NSE:ABCAPITAL"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "NSE:ABCAPITAL"
$ws.Range("C3").Value = "NSE:GOACARBON"
$ws.Range("F3").Value = "NSE:ESCORTS"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "NSE:AHLEAST"
$ws.Range("C4").Value = "NSE:GRANULES"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "NSE:AMNPLST"
$ws.Range("C5").Value = "NSE:GUFICBIO"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "NSE:APOLLO"
$ws.Range("C6").Value = "NSE:ICICIGI"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "NSE:ASALCBR"
$ws.Range("C7").Value = "NSE:JYOTHYLAB"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "NSE:AXISILVER"
$ws.Range("C8").Value = "NSE:MANGLMCEM"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "NSE:BBTC"
$ws.Range("C9").Value = "NSE:ORIENTHOT"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "NSE:BHAGERIA"
$ws.Range("C10").Value = "NSE:PGHL"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "NSE:BHAGYANGR"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "NSE:CDSL"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "NSE:CONFIPET"

# New rows 14-45
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NSE:CREST"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "NSE:DEN"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "NSE:DPWIRES"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "NSE:EMAMIPAP"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "NSE:GRAPHITE"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "NSE:GREAVESCOT"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "NSE:GREENPOWER"
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "NSE:GSEC10YEAR"
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "NSE:GTPL"
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "NSE:HDFCLOWVOL"
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "NSE:HEG"
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "NSE:HEIDELBERG"
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "NSE:HIMATSEIDE"
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "NSE:INDIACEM"
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "NSE:INTLCONV"
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "NSE:KALAMANDIR"
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "NSE:KRBL"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "NSE:KTKBANK"
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "NSE:LLOYDSME"
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "NSE:LTGILTBEES"
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "NSE:MAHSCOOTER"
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "NSE:MIRZAINT"
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "NSE:MOM50"
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "NSE:MONQ50"
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "NSE:MVGJL"
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "NSE:ONMOBILE"
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "NSE:ORISSAMINE"
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "NSE:POWERINDIA"
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "NSE:RADIANTCMS"
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "NSE:RGL"
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "NSE:RITES"
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "NSE:SALONA"
# Copy the bordered/bold/centered style from an existing column-A cell down to the new rows
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A14:A45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
